$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.325.51'
$ws.Range("E2").Value = '  +3.05%  '
$ws.Range("D3").Value = '1.715.80'
$ws.Range("E3").Value = '  +3.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.01'
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4755'
$ws.Range("E7").Value = '  -0.73%  '
$ws.Range("E8").Value = '  +0.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06209'
$ws.Range("E9").Value = '  +0.96%  '
$ws.Range("D10").Value = '1.716.31'
$ws.Range("E10").Value = '  +3.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07048'
$ws.Range("E11").Value = '  -0.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.31'
$ws.Range("E12").Value = '  +3.97%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5883'
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.416'
$ws.Range("E14").Value = '  +1.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '75.98'
$ws.Range("E15").Value = '  +2.25%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '26.334.65'
$ws.Range("E18").Value = '  +3.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006816'
$ws.Range("E19").Value = '  +0.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.54'
$ws.Range("E20").Value = '  +1.18%  '
$ws.Range("D21").Value = '1.934.82'
$ws.Range("E21").Value = '  +2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.540'
$ws.Range("E22").Value = '  +2.61%  '
$ws.Range("E23").Value = '  +1.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.315'
$ws.Range("E24").Value = '  +0.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.18'
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.22'
$ws.Range("E26").Value = '  +1.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '107.95'
$ws.Range("E27").Value = '  +3.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.403'
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.750'
$ws.Range("E29").Value = '  +3.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.991'
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.685'
$ws.Range("E31").Value = '  +0.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07751'
$ws.Range("E32").Value = '  +1.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04432'
$ws.Range("E33").Value = '  +2.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.612'
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9737'
$ws.Range("E35").Value = '  +2.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6187'
$ws.Range("E36").Value = '  +1.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9259'
$ws.Range("E37").Value = '  +9.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '112.19'
$ws.Range("E38").Value = '  +14.54%  '
$ws.Range("E39").Value = '  -7.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.914'
$ws.Range("E40").Value = '  +2.63%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01471'
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.315'
$ws.Range("E43").Value = '  +13.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3806'
$ws.Range("E44").Value = '  +1.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1163'
$ws.Range("E45").Value = '  +3.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.294'
$ws.Range("E46").Value = '  +1.49%  '
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.28'
$ws.Range("E48").Value = '  +2.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.660'
$ws.Range("E49").Value = '  +4.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '50.68'
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3360'
$ws.Range("E51").Value = '  +0.90%  '
